$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 56; existing rows 56-74 shift down to 57-75.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new record.
$ws.Cells.Item(56, 1).Value = 3
$ws.Cells.Item(56, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(56, 3).Value = "Coquimbo"
$ws.Cells.Item(56, 4).Value = 45205
$ws.Cells.Item(56, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(56, 5).Value = 5
$ws.Cells.Item(56, 6).Value = 300000000
$ws.Cells.Item(56, 7).Value = "Espárragos"
$ws.Cells.Item(56, 8).Value = "Verde"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 1200
$ws.Cells.Item(56, 11).Value = 1800
$ws.Cells.Item(56, 12).Value = 1800
$ws.Cells.Item(56, 13).Value = 1800
$ws.Cells.Item(56, 14).Value = "$/kilo"
$ws.Cells.Item(56, 15).Value = "Provincia de Linares"
$ws.Cells.Item(56, 16).Value = 1800
$ws.Cells.Item(56, 17).Value = 1
$ws.Cells.Item(56, 18).Value = "Hortaliza"
